# Update rows 2-4 (existing ECs-sender rows get corrected LR metrics & target-cluster
# strings reshuffled) and append new rows 5-7 for the sCs-sender combinations,
# per Dr Hou's revised Natmi computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cdh1"
$ws.Cells.Item(2, 3).Value = "Itgb7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.098866
$ws.Cells.Item(2, 8).Value = 0.296598
$ws.Cells.Item(2, 9).Value = 0.3026185969870575
$ws.Cells.Item(2, 10).Value = 0.3026185969870575
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 5.864463
$ws.Cells.Item(2, 14).Value = 17.593389
$ws.Cells.Item(2, 15).Value = 0.6069167733108516
$ws.Cells.Item(2, 16).Value = 0.6069167733108515
$ws.Cells.Item(2, 17).Value = 0.5797959989579999
$ws.Cells.Item(2, 18).Value = 5.218163990621999
$ws.Cells.Item(2, 19).Value = 0.1836643024272419
$ws.Cells.Item(2, 20).Value = 0.1836643024272419

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cdh1"
$ws.Cells.Item(3, 3).Value = "Itgb7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.098866
$ws.Cells.Item(3, 8).Value = 0.296598
$ws.Cells.Item(3, 9).Value = 0.3026185969870575
$ws.Cells.Item(3, 10).Value = 0.3026185969870575
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.691504
$ws.Cells.Item(3, 14).Value = 8.074512
$ws.Cells.Item(3, 15).Value = 0.278545354115671
$ws.Cells.Item(3, 16).Value = 0.278545354115671
$ws.Cells.Item(3, 17).Value = 0.266098234464
$ws.Cells.Item(3, 18).Value = 2.394884110176
$ws.Cells.Item(3, 19).Value = 0.08429300425974746
$ws.Cells.Item(3, 20).Value = 0.08429300425974745

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cdh1"
$ws.Cells.Item(4, 3).Value = "Itgb7"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.098866
$ws.Cells.Item(4, 8).Value = 0.296598
$ws.Cells.Item(4, 9).Value = 0.3026185969870575
$ws.Cells.Item(4, 10).Value = 0.3026185969870575
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.106746666666667
$ws.Cells.Item(4, 14).Value = 3.32024
$ws.Cells.Item(4, 15).Value = 0.1145378725734776
$ws.Cells.Item(4, 16).Value = 0.1145378725734776
$ws.Cells.Item(4, 17).Value = 0.1094196159466667
$ws.Cells.Item(4, 18).Value = 0.9847765435199999
$ws.Cells.Item(4, 19).Value = 0.03466129030006816
$ws.Cells.Item(4, 20).Value = 0.03466129030006815

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Cdh1"
$ws.Cells.Item(5, 3).Value = "Itgb7"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.2278356666666667
$ws.Cells.Item(5, 8).Value = 0.6835070000000001
$ws.Cells.Item(5, 9).Value = 0.6973814030129426
$ws.Cells.Item(5, 10).Value = 0.6973814030129426
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.864463
$ws.Cells.Item(5, 14).Value = 17.593389
$ws.Cells.Item(5, 15).Value = 0.6069167733108516
$ws.Cells.Item(5, 16).Value = 0.6069167733108515
$ws.Cells.Item(5, 17).Value = 1.336133837247
$ws.Cells.Item(5, 18).Value = 12.025204535223
$ws.Cells.Item(5, 19).Value = 0.4232524708836097
$ws.Cells.Item(5, 20).Value = 0.4232524708836096

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Cdh1"
$ws.Cells.Item(6, 3).Value = "Itgb7"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.2278356666666667
$ws.Cells.Item(6, 8).Value = 0.6835070000000001
$ws.Cells.Item(6, 9).Value = 0.6973814030129426
$ws.Cells.Item(6, 10).Value = 0.6973814030129426
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.691504
$ws.Cells.Item(6, 14).Value = 8.074512
$ws.Cells.Item(6, 15).Value = 0.278545354115671
$ws.Cells.Item(6, 16).Value = 0.278545354115671
$ws.Cells.Item(6, 17).Value = 0.6132206081760001
$ws.Cells.Item(6, 18).Value = 5.518985473584001
$ws.Cells.Item(6, 19).Value = 0.1942523498559236
$ws.Cells.Item(6, 20).Value = 0.1942523498559236

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Cdh1"
$ws.Cells.Item(7, 3).Value = "Itgb7"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.2278356666666667
$ws.Cells.Item(7, 8).Value = 0.6835070000000001
$ws.Cells.Item(7, 9).Value = 0.6973814030129426
$ws.Cells.Item(7, 10).Value = 0.6973814030129426
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.106746666666667
$ws.Cells.Item(7, 14).Value = 3.32024
$ws.Cells.Item(7, 15).Value = 0.1145378725734776
$ws.Cells.Item(7, 16).Value = 0.1145378725734776
$ws.Cells.Item(7, 17).Value = 0.2521563646311111
$ws.Cells.Item(7, 18).Value = 2.26940728168
$ws.Cells.Item(7, 19).Value = 0.07987658227340942
$ws.Cells.Item(7, 20).Value = 0.07987658227340942
